# Update the "Lương" sheet with the corrected values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 19
$ws.Range("B3").Value = 665000
$ws.Range("B12").Value = 2714285.714285714
$ws.Range("B29").Value = 3559285.714285714
$ws.Range("B31").Value = 3559285.714285714
